# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.631.04'
$ws.Range('E2').Value = '  +4.63%  '
$ws.Range('D3').Value = '3.493.34'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +8.89%  '
$ws.Range('D9').Value = '3.490.23'
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('E10').Value = '  +7.51%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.439'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.27%  '
$ws.Range('D13').Value = '4.098.72'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.11%  '
$ws.Range('E16').Value = '  +4.30%  '
$ws.Range('D17').Value = '66.619.84'
$ws.Range('E17').Value = '  +4.53%  '
$ws.Range('D18').Value = '3.496.14'
$ws.Range('E18').Value = '  +2.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '392.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.86%  '
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  +4.76%  '
$ws.Range('E26').Value = '  +6.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.14%  '
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  +5.12%  '
$ws.Range('E31').Value = '  +7.11%  '
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('E34').Value = '  +5.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.52'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.884'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.32%  '
$ws.Range('E39').Value = '  +6.89%  '
$ws.Range('E40').Value = '  +7.60%  '
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.85%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '27.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').Value = '2.781.19'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.68%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '352.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.29%  '
$ws.Range('E50').Value = '  +6.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +14.44%  '
